$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 308-313 with new values
$ws.Range("B308").Value = 43
$ws.Range("C308").Value = 151
$ws.Range("D308").Value = 2
$ws.Range("E308").Value = 19

$ws.Range("B309").Value = 26
$ws.Range("C309").Value = 86
$ws.Range("D309").Value = 1
$ws.Range("E309").Value = 22

$ws.Range("B310").Value = 27
$ws.Range("C310").Value = 90
$ws.Range("D310").Value = 1
$ws.Range("E310").Value = 22

$ws.Range("B311").Value = 54
$ws.Range("C311").Value = 168
$ws.Range("D311").Value = 3
$ws.Range("E311").Value = 154

$ws.Range("B312").Value = 20
$ws.Range("C312").Value = 63
$ws.Range("D312").Value = 1
$ws.Range("E312").Value = 102

$ws.Range("B313").Value = 15
$ws.Range("C313").Value = 78
$ws.Range("D313").Value = 1
$ws.Range("E313").Value = 102

# Add new row 314 - force column A to stay plain text (not auto-converted to a date serial)
$ws.Range("A314").NumberFormat = "@"
$ws.Range("A314").Value = "11.01.2021"
$ws.Range("A314").ClearFormats()

$ws.Range("B314").Value = 6
$ws.Range("C314").Value = 43
$ws.Range("D314").Value = 3
$ws.Range("E314").Value = 111
